$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.953047333333333
$ws.Range("H2").Value = 14.859142
$ws.Range("I2").Value = 0.7703204220313993
$ws.Range("J2").Value = 0.7703204220313993
$ws.Range("O2").Value = 0.01909956851648506
$ws.Range("P2").Value = 0.01909956851648506
$ws.Range("Q2").Value = 0.07167884999222222
$ws.Range("R2").Value = 0.64510964993
$ws.Range("S2").Value = 0.0147127876802364
$ws.Range("T2").Value = 0.0147127876802364

$ws.Range("G3").Value = 4.953047333333333
$ws.Range("H3").Value = 14.859142
$ws.Range("I3").Value = 0.7703204220313993
$ws.Range("J3").Value = 0.7703204220313993
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7432243333333334
$ws.Range("N3").Value = 2.229673
$ws.Range("O3").Value = 0.9809004314835149
$ws.Range("P3").Value = 0.9809004314835149
$ws.Range("Q3").Value = 3.681225302285111
$ws.Range("R3").Value = 33.131027720566
$ws.Range("S3").Value = 0.7556076343511628
$ws.Range("T3").Value = 0.7556076343511628

$ws.Range("I4").Value = 0.135969508894967
$ws.Range("J4").Value = 0.135969508894967
$ws.Range("O4").Value = 0.01909956851648506
$ws.Range("P4").Value = 0.01909956851648506
$ws.Range("S4").Value = 0.002596958951292248
$ws.Range("T4").Value = 0.002596958951292248

$ws.Range("I5").Value = 0.135969508894967
$ws.Range("J5").Value = 0.135969508894967
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7432243333333334
$ws.Range("N5").Value = 2.229673
$ws.Range("O5").Value = 0.9809004314835149
$ws.Range("P5").Value = 0.9809004314835149
$ws.Range("Q5").Value = 0.6497742785573334
$ws.Range("R5").Value = 5.847968507016
$ws.Range("S5").Value = 0.1333725499436748
$ws.Range("T5").Value = 0.1333725499436748

$ws.Range("G6").Value = 0.5382536666666667
$ws.Range("H6").Value = 1.614761
$ws.Range("I6").Value = 0.08371165542397027
$ws.Range("J6").Value = 0.08371165542397027
$ws.Range("O6").Value = 0.01909956851648506
$ws.Range("P6").Value = 0.01909956851648506
$ws.Range("Q6").Value = 0.007789427646111112
$ws.Range("R6").Value = 0.07010484881500001
$ws.Range("S6").Value = 0.001598856498398509
$ws.Range("T6").Value = 0.001598856498398509

$ws.Range("G7").Value = 0.5382536666666667
$ws.Range("H7").Value = 1.614761
$ws.Range("I7").Value = 0.08371165542397027
$ws.Range("J7").Value = 0.08371165542397027
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7432243333333334
$ws.Range("N7").Value = 2.229673
$ws.Range("O7").Value = 0.9809004314835149
$ws.Range("P7").Value = 0.9809004314835149
$ws.Range("Q7").Value = 0.4000432225725556
$ws.Range("R7").Value = 3.600389003153
$ws.Range("S7").Value = 0.08211279892557176
$ws.Range("T7").Value = 0.08211279892557176

$ws.Range("G8").Value = 0.06428833333333334
$ws.Range("H8").Value = 0.192865
$ws.Range("I8").Value = 0.009998413649663342
$ws.Range("J8").Value = 0.009998413649663342
$ws.Range("O8").Value = 0.01909956851648506
$ws.Range("P8").Value = 0.01909956851648506
$ws.Range("Q8").Value = 0.0009303593305555557
$ws.Range("R8").Value = 0.008373233975000002
$ws.Range("S8").Value = 0.0001909653865579045
$ws.Range("T8").Value = 0.0001909653865579045

$ws.Range("G9").Value = 0.06428833333333334
$ws.Range("H9").Value = 0.192865
$ws.Range("I9").Value = 0.009998413649663342
$ws.Range("J9").Value = 0.009998413649663342
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7432243333333334
$ws.Range("N9").Value = 2.229673
$ws.Range("O9").Value = 0.9809004314835149
$ws.Range("P9").Value = 0.9809004314835149
$ws.Range("Q9").Value = 0.04778065368277779
$ws.Range("R9").Value = 0.430025883145
$ws.Range("S9").Value = 0.009807448263105437
$ws.Range("T9").Value = 0.009807448263105437
